# Auto-generated edit script: update cryptos.xlsx price/volume/coin data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell 'D2' '30.059.14'
Set-TextCell 'E2' '  -2.49%  '
Set-TextCell 'D3' '1.857.24'
Set-TextCell 'E3' '  -3.60%  '
Set-TextCell 'D4' '1.002'
Set-TextCell 'E4' '  +0.25%  '
Set-TextCell 'D5' '233.39'
Set-TextCell 'E5' '  -3.46%  '
Set-TextCell 'E6' '  +0.18%  '
Set-TextCell 'D7' '0.4639'
Set-TextCell 'E7' '  -3.14%  '
Set-TextCell 'D8' '0.2803'
Set-TextCell 'E8' '  -3.01%  '
Set-TextCell 'D9' '0.06531'
Set-TextCell 'E9' '  -3.78%  '
Set-TextCell 'D10' '19.47'
Set-TextCell 'E10' '  -1.09%  '
Set-TextCell 'D11' '0.07806'
Set-TextCell 'E11' '  +0.15%  '
Set-TextCell 'D12' '96.45'
Set-TextCell 'E12' '  -7.60%  '
Set-TextCell 'D13' '1.863.51'
Set-TextCell 'E13' '  -3.85%  '
Set-TextCell 'D14' '5.092'
Set-TextCell 'E14' '  -3.55%  '
Set-TextCell 'D15' '0.6614'
Set-TextCell 'E15' '  -3.40%  '
Set-TextCell 'D16' '279.16'
Set-TextCell 'E16' '  -4.68%  '
Set-TextCell 'D17' '30.109.80'
Set-TextCell 'E17' '  -2.34%  '
Set-TextCell 'D18' '1.001'
Set-TextCell 'E18' '  +0.10%  '
Set-TextCell 'D19' '5.440'
Set-TextCell 'E19' '  -1.50%  '
Set-TextCell 'D20' '12.55'
Set-TextCell 'E20' '  -2.57%  '
Set-TextCell 'D21' '2.108.78'
Set-TextCell 'E21' '  -3.47%  '
Set-TextCell 'B22' 'BinanceUSD'
Set-TextCell 'C22' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell 'D22' '1.003'
Set-TextCell 'E22' '  +0.31%  '
Set-TextCell 'B23' 'ShibaInu'
Set-TextCell 'C23' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell 'D23' '0.000007213'
Set-TextCell 'E23' '  -5.03%  '
Set-TextCell 'D24' '6.114'
Set-TextCell 'E24' '  -4.32%  '
Set-TextCell 'D25' '9.303'
Set-TextCell 'E25' '  -2.58%  '
Set-TextCell 'D26' '167.03'
Set-TextCell 'E26' '  -0.54%  '
Set-TextCell 'D27' '18.83'
Set-TextCell 'E27' '  -4.81%  '
Set-TextCell 'D28' '1.902'
Set-TextCell 'E28' '  -10.22%  '
Set-TextCell 'D29' '1.339'
Set-TextCell 'E29' '  -3.74%  '
Set-TextCell 'D30' '0.09532'
Set-TextCell 'D31' '4.416'
Set-TextCell 'E31' '  -4.08%  '
Set-TextCell 'D32' '1.468'
Set-TextCell 'E32' '  -4.02%  '
Set-TextCell 'D33' '4.093'
Set-TextCell 'E33' '  -5.45%  '
Set-TextCell 'D34' '0.04627'
Set-TextCell 'E34' '  -3.87%  '
Set-TextCell 'D35' '1.094'
Set-TextCell 'E35' '  -2.83%  '
Set-TextCell 'D36' '0.6966'
Set-TextCell 'E36' '  -5.40%  '
Set-TextCell 'D37' '2.708'
Set-TextCell 'E37' '  -0.39%  '
Set-TextCell 'D38' '0.01833'
Set-TextCell 'E38' '  -5.75%  '
Set-TextCell 'D39' '6.307'
Set-TextCell 'E39' '  -1.70%  '
Set-TextCell 'D40' '2.511'
Set-TextCell 'E40' '  -4.73%  '
Set-TextCell 'D41' '71.91'
Set-TextCell 'E41' '  -4.51%  '
Set-TextCell 'D42' '0.8527'
Set-TextCell 'E42' '  -1.88%  '
Set-TextCell 'D43' '1.910'
Set-TextCell 'E43' '  -5.22%  '
Set-TextCell 'B44' 'PaxDollar'
Set-TextCell 'C44' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell 'D44' '1.001'
Set-TextCell 'E44' '  +0.12%  '
Set-TextCell 'B45' 'Quant'
Set-TextCell 'C45' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 'D45' '103.65'
Set-TextCell 'E45' '  -2.07%  '
Set-TextCell 'D46' '0.4124'
Set-TextCell 'E46' '  -5.17%  '
Set-TextCell 'D47' '989.73'
Set-TextCell 'E47' '  -0.62%  '
Set-TextCell 'D48' '7.163'
Set-TextCell 'E48' '  -5.05%  '
Set-TextCell 'D49' '9.236'
Set-TextCell 'E49' '  +2.51%  '
Set-TextCell 'D50' '33.88'
Set-TextCell 'E50' '  -3.05%  '
Set-TextCell 'D51' '0.1135'
Set-TextCell 'E51' '  -6.35%  '
